$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for column C (rows 2-7)
$ws.Range("C2").Value = 0.4278141
$ws.Range("C3").Value = 0.6195014999999999
$ws.Range("C4").Value = 0.8120718
$ws.Range("C5").Value = 1.0068984
$ws.Range("C6").Value = 1.2043737
$ws.Range("C7").Value = 1.4004756

# New values for column E (rows 2-7) - all identical
$eVal = 0.01001948728320346
$ws.Range("E2").Value = $eVal
$ws.Range("E3").Value = $eVal
$ws.Range("E4").Value = $eVal
$ws.Range("E5").Value = $eVal
$ws.Range("E6").Value = $eVal
$ws.Range("E7").Value = $eVal
